# slides-118-BRSKI-AE.pptx: adapted layout on 2nd slide
#
# The second slide's body placeholder ("CustomShape 2") had the font size
# of its first two bullets (and the second bullet's nested sub-bullets,
# through the closing quote after "...brski-registrar-cmp...") reduced
# from 22pt to 20pt so the text fits better on the slide.
#
# The affected runs form one contiguous block of text, starting at the
# very first character of the shape and running through the closing
# curly quote that ends the "brski-registrar-cmp" bullet, so we can
# resize them in a single Characters() call instead of touching each
# run individually.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# The body text box is "CustomShape 2" (the 2nd shape on the slide); locate
# it by name so the script still works if shape ordering ever shifts.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "CustomShape 2") {
        $shape = $candidate
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(2)
}

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Characters(Start, Length) -- 1-based, inclusive range covering:
#   "SECDIR Last Call Review completed by Barry Leiba on Nov 4th; status: Ready"
#   "Alignment during IETF 118 regarding discovery of registrars with enhanced feature sets:"
#   "Accepted as general problem (arrow) general solution addressed in new draft: [BRSKI-Discovery]"
#   "Absence of general solution handled in BRSKI-AE by specific service name "brski-registrar-cmp""
$target = $tr.Characters(1, 346)
$target.Font.Size = 20
